{"js": "// Insert the group's data as a run of text at the very start of the\n// document's (only) paragraph, ahead of the existing _GoBack bookmark.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// The target document splits this into three <w:r> runs around\n// proofing-error markers Word's spell-checker added for \"Selman\" while\n// it was typed; those markers are a transient UI artifact (not part of\n// the Office.js/COM object models), so a single insertText() call that\n// lands the same visible text ahead of the existing _GoBack bookmark\n// reproduces the same reading-order content and run layout that matters.\nfirstParagraph.insertText(\n  \"Carlos Andr\u00e9s Selman David \u2013 2018325 \u2013 IN5BM \u2013 Dise\u00f1o de la Aplicaci\u00f3n\",\n  Word.InsertLocation.start\n);\n\nawait context.sync();\n", "ps1": "# Insert the group's data as text at the very start of the document's\n# (only) paragraph, ahead of the existing _GoBack bookmark.\n$d = $word.ActiveDocument\n\n$para = $d.Paragraphs(1)\n$rng = $para.Range\n$rng.Collapse(1)  # wdCollapseStart\n$rng.InsertBefore(\"Carlos Andr\u00e9s Selman David \u2013 2018325 \u2013 IN5BM \u2013 Dise\u00f1o de la Aplicaci\u00f3n\")\n"}
